$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "이중 진자와 카오스 이론"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/06/09/double_pendulum_and_chaos_theory.html"

$ws.Range("D9").Value = "온라인 교육 – “산지직송” 대학"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/online-middle-man/#utm_source=rss&utm_medium=rss&utm_campaign=online-middle-man"

$ws.Range("D28").Value = "[강화학습] Dynamic programming"
$ws.Range("E28").Value = "https://ropiens.tistory.com/130"

$ws.Range("D46").Value = "[질병관리청] 2021년 06월, 생물정보학(Bioinformatics 채용), 공무직(연구원) 채용"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/401"
